# INF 301 - Module 2 - Git
#
# Adds a new slide at the end of the deck: a small "git checkout" diagram
# (Local repository / Working directory boxes on branch "A", linked by an
# arrow, labelled "git checkout"). This diagram is a trimmed-down copy of
# the big Git-concepts diagram that already lives on slide 1, so the new
# slide is produced the same way a human author would do it in the UI:
# duplicate slide 1, move the duplicate to the end of the deck, then
# delete every shape that isn't part of this particular sub-diagram.

$p = $ppt.ActivePresentation

# Duplicate the first slide (it owns the full Git-concepts drawing that
# every other "sub-diagram" slide in this deck is derived from).
$source = $p.Slides.Item(1)
$dup = $source.Duplicate()
$newSlide = $dup.Item(1)

# Move the new slide to the end of the presentation.
$newSlide.MoveTo($p.Slides.Count)

# Only these four shapes belong to the "git checkout" sub-diagram:
#   id 8  -> "Local repository, branch A" rounded rectangle
#   id 11 -> "Working directory, branch A" rounded rectangle
#   id 46 -> arrow connector between the two boxes
#   id 76 -> "git checkout" label textbox
# Everything else that came along with the duplicated slide gets removed.
$keepIds = @(8, 11, 46, 76)

for ($i = $newSlide.Shapes.Count; $i -ge 1; $i--) {
    $shape = $newSlide.Shapes.Item($i)
    if ($keepIds -notcontains $shape.Id) {
        $shape.Delete()
    }
}
